# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-8.
# Column F (Win) is left untouched, matching the diff.
$data = @{
    2 = @{ B = 0.0008583669626518464; C = 0.3127903958511391;  D = 3.900430680208489;  E = 8.660232485948974;  G = 12.87431192897125 }
    3 = @{ B = 0.3048080303191223;    C = 0.3127903958511391;  D = 26.21740644021617;  E = 0.496779210170732;  G = 27.33178407655716 }
    4 = @{ B = 3.230985683306322;     C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    5 = @{ B = 0.003994804209775715;  C = 0.3127903958511391;  D = 0.8054896365839992; E = 8.660232485948974;  G = 9.782507322593888 }
    6 = @{ B = 0.6753301551942219;    C = 114.8270160096505;   D = 3.900430680208489;  E = 8.660232485948974;  G = 128.0630093310022 }
    7 = @{ B = 0.6753301551942219;    C = 1.667794583268128;   D = 337.1190423067083;  E = 8.660232485948974;  G = 348.1223995311196 }
    8 = @{ B = 1.459612070389937;     C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;  G = 7.524616544037286 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
